$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 30011.5
$ws.Range("J108").Value = 30011.5
$ws.Range("L108").Value = 30011.5
$ws.Range("N108").Value = -37691.5
$ws.Range("H109").Value = 33301.668
$ws.Range("J109").Value = 33301.668
$ws.Range("L109").Value = 33301.668
$ws.Range("N109").Value = -36075.668
$ws.Range("H116").Value = 6045.55
$ws.Range("I116").Value = 4100.5
$ws.Range("J116").Value = 6879.143
$ws.Range("K116").Value = 4100.5
$ws.Range("L116").Value = 6879.143
$ws.Range("M116").Value = -658.5
$ws.Range("N116").Value = -13763.143
$ws.Range("H117").Value = 48538
$ws.Range("J117").Value = 48538
$ws.Range("L117").Value = 48538
$ws.Range("N117").Value = -57716
$ws.Range("H120").Value = 49726
$ws.Range("J120").Value = 49726
$ws.Range("L120").Value = 49726
$ws.Range("N120").Value = -59402
$ws.Range("H121").Value = 2565296.8
$ws.Range("J121").Value = 3705297
$ws.Range("L121").Value = 11115891
$ws.Range("N121").Value = -11119385
$ws.Range("H131").Value = 3369.5715
$ws.Range("I131").Value = 2746.25
$ws.Range("J131").Value = 3618.9
$ws.Range("K131").Value = 8238.75
$ws.Range("L131").Value = 10856.7
$ws.Range("M131").Value = -3198.75
$ws.Range("N131").Value = -20936.7
$ws.Range("H132").Value = 29125.766
$ws.Range("I132").Value = 3949.25
$ws.Range("J132").Value = 146616.17
$ws.Range("K132").Value = 11847.75
$ws.Range("L132").Value = 439848.51
$ws.Range("M132").Value = -9317.75
$ws.Range("N132").Value = -444908.51
$ws.Range("H137").Value = 4072.9792
$ws.Range("I137").Value = 1350
$ws.Range("J137").Value = 4389.6045
$ws.Range("K137").Value = 4050
$ws.Range("L137").Value = 13168.8135
$ws.Range("M137").Value = -1500
$ws.Range("N137").Value = -18268.8135
$ws.Range("H138").Value = 2027.8677
$ws.Range("I138").Value = 2095.818
$ws.Range("J138").Value = 1995.3695
$ws.Range("K138").Value = 6287.454000000001
$ws.Range("L138").Value = 5986.1085
$ws.Range("M138").Value = -1147.454000000001
$ws.Range("N138").Value = -16266.1085

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33792.66
$ws.Range("I32").Value = 35829.047
$ws.Range("K32").Value = 35829.047
$ws.Range("M32").Value = -35542.047
$ws.Range("H61").Value = 3862.3872
$ws.Range("I61").Value = 2198.4443
$ws.Range("J61").Value = 4543.091
$ws.Range("K61").Value = 2198.4443
$ws.Range("L61").Value = 4543.091
$ws.Range("M61").Value = -1986.4443
$ws.Range("N61").Value = -4967.091
$ws.Range("H74").Value = 1308.5927
$ws.Range("I74").Value = 629.7646999999999
$ws.Range("J74").Value = 2462.6
$ws.Range("K74").Value = 629.7646999999999
$ws.Range("L74").Value = 2462.6
$ws.Range("M74").Value = 244.2353000000001
$ws.Range("N74").Value = -4210.6
$ws.Range("H77").Value = 1308.5927
$ws.Range("I77").Value = 629.7646999999999
$ws.Range("J77").Value = 2462.6
$ws.Range("K77").Value = 3148.8235
$ws.Range("L77").Value = 12313
$ws.Range("M77").Value = 1219.1765
$ws.Range("N77").Value = -21049
$ws.Range("H107").Value = 38816
$ws.Range("J107").Value = 38816
$ws.Range("L107").Value = 38816
$ws.Range("N107").Value = -46496
$ws.Range("H132").Value = 18520528
$ws.Range("J132").Value = 2877.9
$ws.Range("L132").Value = 8633.700000000001
$ws.Range("N132").Value = -13693.7
$ws.Range("H134").Value = 52214.285
$ws.Range("J134").Value = 52214.285
$ws.Range("L134").Value = 52214.285
$ws.Range("N134").Value = -62354.285
$ws.Range("H136").Value = 3862.3872
$ws.Range("I136").Value = 2198.4443
$ws.Range("J136").Value = 4543.091
$ws.Range("K136").Value = 6595.3329
$ws.Range("L136").Value = 13629.273
$ws.Range("M136").Value = -4045.3329
$ws.Range("N136").Value = -18729.273

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H132").Value = 47197.145
$ws.Range("J132").Value = 47197.145
$ws.Range("L132").Value = 47197.145
$ws.Range("N132").Value = -57317.145
$ws.Range("H134").Value = 4803.8306
$ws.Range("I134").Value = 4144.4
$ws.Range("J134").Value = 4864.8887
$ws.Range("K134").Value = 12433.2
$ws.Range("L134").Value = 14594.6661
$ws.Range("M134").Value = -9898.199999999999
$ws.Range("N134").Value = -19664.6661
$ws.Range("H140").Value = 37799
$ws.Range("J140").Value = 37799
$ws.Range("L140").Value = 37799
$ws.Range("N140").Value = -48159

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 331766.84
$ws.Range("I31").Value = 5867
$ws.Range("K31").Value = 5867
$ws.Range("M31").Value = -5572
$ws.Range("H34").Value = 331766.84
$ws.Range("I34").Value = 5867
$ws.Range("K34").Value = 5867
$ws.Range("M34").Value = -5665
$ws.Range("H86").Value = 5361
$ws.Range("I86").Value = 5361
$ws.Range("K86").Value = 5361
$ws.Range("M86").Value = -4238
$ws.Range("H89").Value = 5361
$ws.Range("I89").Value = 5361
$ws.Range("K89").Value = 26805
$ws.Range("M89").Value = -21189
$ws.Range("H132").Value = 63658
$ws.Range("I132").Value = 2424.5715
$ws.Range("J132").Value = 158910
$ws.Range("K132").Value = 7273.7145
$ws.Range("L132").Value = 476730
$ws.Range("M132").Value = -4743.7145
$ws.Range("N132").Value = -481790
$ws.Range("H134").Value = 1402552.2
$ws.Range("I134").Value = 1522.2
$ws.Range("J134").Value = 2803582.2
$ws.Range("K134").Value = 4566.6
$ws.Range("L134").Value = 8410746.600000001
$ws.Range("M134").Value = -2031.6
$ws.Range("N134").Value = -8415816.600000001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 251450.67
$ws.Range("I44").Value = 600800
$ws.Range("J44").Value = 1915.4286
$ws.Range("K44").Value = 1802400
$ws.Range("L44").Value = 5746.2858
$ws.Range("M44").Value = -1802002
$ws.Range("N44").Value = -6542.2858
$ws.Range("H132").Value = 3138.625
$ws.Range("I132").Value = 668.8182
$ws.Range("J132").Value = 5228.4614
$ws.Range("K132").Value = 6019.3638
$ws.Range("L132").Value = 47056.1526
$ws.Range("M132").Value = -3489.3638
$ws.Range("N132").Value = -52116.1526
$ws.Range("H137").Value = 50010170
$ws.Range("I137").Value = 3712.8572
$ws.Range("J137").Value = 76936720
$ws.Range("K137").Value = 11138.5716
$ws.Range("L137").Value = 230810160
$ws.Range("M137").Value = -6038.571599999999
$ws.Range("N137").Value = -230820360

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 46650.5
$ws.Range("J110").Value = 46650.5
$ws.Range("L110").Value = 46650.5
$ws.Range("N110").Value = -54830.5
$ws.Range("H132").Value = 8063.091
$ws.Range("I132").Value = 1200
$ws.Range("K132").Value = 3600
$ws.Range("M132").Value = -1070

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 43946.5
$ws.Range("J111").Value = 43946.5
$ws.Range("L111").Value = 43946.5
$ws.Range("N111").Value = -52126.5
$ws.Range("H132").Value = 3772.3667
$ws.Range("I132").Value = 3009.4211
$ws.Range("K132").Value = 9028.263300000001
$ws.Range("M132").Value = -6498.263300000001

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 36818
$ws.Range("J108").Value = 36818
$ws.Range("L108").Value = 36818
$ws.Range("N108").Value = -44498
$ws.Range("H133").Value = 81991.39999999999
$ws.Range("J133").Value = 81991.39999999999
$ws.Range("L133").Value = 81991.39999999999
$ws.Range("N133").Value = -92111.39999999999
